# The deck originally carried a custom "Integral" colour theme on
# ppt/theme/theme1.xml (the presentation / slide-master theme) and the
# generic default "Office Theme" palette on ppt/theme/theme2.xml (the
# notes-master theme). The edit swaps the two palettes so the slide
# master now uses the stock Office colours.
#
# PowerPoint's COM model exposes a theme's twelve colour slots through
# Theme.ThemeColorScheme (1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1-6,
# 11=hlink, 12=folHlink) as RGB (VBA-style 0xBBGGRR) integers, so drive
# the swap through that object rather than touching the package parts
# directly.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Standard Office theme colours (what ppt/theme/theme2.xml already held).
$colors.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
